# Generate Report for Handback
# Refreshes the handback-status report: the first source file's generated
# UUID/hash/timestamps move forward, and the second source file gets a
# brand-new UUID ("ffff778366a4-ef15-40b7-a639-bcc0e5053014.md").

$wb = $excel.ActiveWorkbook

$oldUuid1 = "39b05a38-cb78-449e-9597-5c512e044c2c"
$oldUuid2 = "b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce"
$newUuid1 = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb"
$newUuid2 = "ffff778366a4-ef15-40b7-a639-bcc0e5053014"

$newHash1 = "4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf"

$newFile1 = "$newUuid1.md"
$newFile2 = "$newUuid2.md"
$newFile1Disp = "e2e\$newFile1"
$newFile2Disp = "e2e\$newFile2"

$newLatestHoDate = "2016-09-06 07:14:00"

$newZhXlf1 = "$newUuid1.$newHash1.zh-cn.xlf"
$newZhHandoff = "2016-09-06 07:13:55"
$newZhHandback = "2016-09-06 07:14:20"

$newDeXlf1 = "$newUuid1.$newHash1.de-de.xlf"
$newDeHandback = "2016-09-06 07:14:28"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newFile1
$ws1.Range("B2").Value = $newFile1Disp
$ws1.Range("G2").Value = $newLatestHoDate

$ws1.Range("A3").Value = $newFile2
$ws1.Range("B3").Value = $newFile2Disp
$ws1.Range("G3").Value = $newLatestHoDate

# Rebuild the hyperlinks on column B, keeping the same target addresses but
# refreshing the display text to match the new file names.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid1.md", "", "", $newFile1Disp)
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid2.md", "", "", $newFile2Disp)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newFile1
$ws2.Range("G2").Value = $newZhXlf1
$ws2.Range("H2").Value = $newZhHandoff
$ws2.Range("I2").Value = $newFile1
$ws2.Range("J2").Value = $newZhXlf1
$ws2.Range("K2").Value = $newZhHandback

$ws2.Range("A3").Value = $newFile2
$ws2.Range("G3").Value = $newZhXlf1
$ws2.Range("H3").Value = $newZhHandoff
$ws2.Range("I3").Value = $newFile2
$ws2.Range("J3").Value = $newZhXlf1
$ws2.Range("K3").Value = $newZhHandback

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid1.md", "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/56a57331e17b61ce845e60cdff58bbddc47edead/e2e/$oldUuid1.md", "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid2.md", "", "", $newFile2)
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/56a57331e17b61ce845e60cdff58bbddc47edead/e2e/$oldUuid2.md", "", "", $newFile2)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newFile1
$ws3.Range("G2").Value = $newDeXlf1
$ws3.Range("H2").Value = $newLatestHoDate
$ws3.Range("I2").Value = $newFile1
$ws3.Range("J2").Value = $newDeXlf1
$ws3.Range("K2").Value = $newDeHandback

$ws3.Range("A3").Value = $newFile2
$ws3.Range("G3").Value = $newDeXlf1
$ws3.Range("H3").Value = $newLatestHoDate
$ws3.Range("I3").Value = $newFile2
$ws3.Range("J3").Value = $newDeXlf1
$ws3.Range("K3").Value = $newDeHandback

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid1.md", "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2a95a885c32c3eb4eb821a6807913ca09502c90d/e2e/$oldUuid1.md", "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c727d2f8499be35219865b79759b978f7ac3c4a/e2e/$oldUuid2.md", "", "", $newFile2)
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2a95a885c32c3eb4eb821a6807913ca09502c90d/e2e/$oldUuid2.md", "", "", $newFile2)

Write-Output "Handback status report regenerated."
